$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) About sheet: bump the "last updated" date in C1
#    45366 (2024-03-15) -> 45379 (2024-03-28)
# ------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = 45379

# ------------------------------------------------------------------
# 2) RAF-capacity sheet: raise the capacity-credit multiplier for the
#    two hydrogen technologies (hydrogen combustion turbine / hydrogen
#    combined cycle) from 0.3 up to 1
# ------------------------------------------------------------------
$wsCapacity = $wb.Worksheets.Item("RAF-capacity")
$wsCapacity.Range("B24").Value = 1
$wsCapacity.Range("B25").Value = 1

# widen column A a touch on the RAF-capacity sheet
$wsCapacity.Columns("A").ColumnWidth = 28.1

# ------------------------------------------------------------------
# 3) View/navigation state: RAF-capacity becomes the active, selected
#    tab (previously RAF-generation was selected/active); zoom in to
#    80% and leave the selection on B25.
# ------------------------------------------------------------------
$wsCapacity.Activate()
$wsCapacity.Range("B25").Select()
$excel.ActiveWindow.Zoom = 80
